{"js": "// Add a new \"SOPickList - MASTER - DO NOT MODIFY\" bullet right after the\n// existing \"... - Quote Template - MASTER - DO NOT MODIFY\" bullet, matching\n// the same List Paragraph / numbering / blue font formatting used by its\n// sibling bullets in that list.\n\nconst body = context.document.body;\nbody.paragraphs.load('items,text');\nawait context.sync();\n\n// Find the paragraph that ends the \"Quote Template - MASTER - DO NOT MODIFY\" line.\nconst paragraphs = body.paragraphs.items;\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.length; i++) {\n  const text = paragraphs[i].text;\n  if (text.indexOf('Quote Template - MASTER - DO NOT MODIFY') !== -1) {\n    targetParagraph = paragraphs[i];\n    break;\n  }\n}\n\nif (!targetParagraph) {\n  throw new Error('Could not locate the \"Quote Template - MASTER - DO NOT MODIFY\" paragraph.');\n}\n\n// Inserting \"after\" an existing list paragraph clones its paragraph\n// properties (style, numbering, spacing) and run formatting automatically.\nconst newParagraph = targetParagraph.insertParagraph(\n  'SOPickList - MASTER - DO NOT MODIFY',\n  Word.InsertLocation.after\n);\nnewParagraph.font.color = '#0070C0';\nawait context.sync();\n\n// Re-load so the range we grab below reflects the just-inserted text.\nbody.paragraphs.load('items,text');\nawait context.sync();\n\nconst refreshedParagraphs = body.paragraphs.items;\nlet newParagraphFresh = null;\nfor (let i = 0; i < refreshedParagraphs.length; i++) {\n  if (refreshedParagraphs[i].text === 'SOPickList - MASTER - DO NOT MODIFY') {\n    newParagraphFresh = refreshedParagraphs[i];\n  }\n}\n\nif (newParagraphFresh) {\n  const endRange = newParagraphFresh.getRange('End');\n  endRange.insertBookmark('_GoBack');\n  await context.sync();\n}\n", "ps1": "# Add a new \"SOPickList - MASTER - DO NOT MODIFY\" bullet right after the\n# existing \"... - Quote Template - MASTER - DO NOT MODIFY\" bullet, matching\n# the same List Paragraph / numbering / blue font formatting used by its\n# sibling bullets in that list, and tag the end of the new bullet with the\n# \"_GoBack\" bookmark Word leaves at the most recent edit location.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that ends the \"Quote Template - MASTER - DO NOT MODIFY\" line.\n$paragraphs = $d.Paragraphs\n$paragraphCount = $paragraphs.Count\n$targetIndex = -1\nfor ($i = 1; $i -le $paragraphCount; $i++) {\n    if ($paragraphs.Item($i).Range.Text -like \"*Quote Template - MASTER - DO NOT MODIFY*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw 'Could not locate the \"Quote Template - MASTER - DO NOT MODIFY\" paragraph.'\n}\n\n$targetParagraph = $paragraphs.Item($targetIndex)\n\n# Collapse to the end of that paragraph (after its text, before the pilcrow)\n# and insert a brand-new paragraph mark after it; this clones the paragraph's\n# style/numbering/spacing formatting automatically, same as the sibling bullets.\n$endOfTarget = $targetParagraph.Range\n$endOfTarget.Collapse(0)   # wdCollapseEnd\n$endOfTarget.InsertParagraphAfter()\n\n# Re-fetch the paragraph collection so we get a live reference to the new,\n# still-empty paragraph that was just created.\n$paragraphs = $d.Paragraphs\n$newParagraph = $paragraphs.Item($targetIndex + 1)\n$newRange = $newParagraph.Range\n\n# Set the text (with a one-character sentinel tacked on) and the blue font\n# color (RGB 0,112,192 -> 0x0070C0, stored by Word as a BGR-packed long).\n$bodyText = 'SOPickList - MASTER - DO NOT MODIFY'\n$sentinel = 'X'\n$newRange.Text = $bodyText + $sentinel\n\n$red = 0x00\n$green = 0x70\n$blue = 0xC0\n$newRange.Font.Color = ($blue * 65536) + ($green * 256) + $red\n\n# Re-fetch again after the Text assignment so positions are current, then\n# drop a \"_GoBack\" bookmark right after the real text but before the\n# sentinel character (an interior, unambiguous position), and finally\n# delete the sentinel so the bookmark is left sitting at the true end of\n# the paragraph's text, immediately before the paragraph mark.\n$paragraphs = $d.Paragraphs\n$newParagraph = $paragraphs.Item($targetIndex + 1)\n$paragraphStart = $newParagraph.Range.Start\n$bodyLength = $bodyText.Length\n\n$bookmarkRange = $d.Range($paragraphStart + $bodyLength, $paragraphStart + $bodyLength)\n$d.Bookmarks.Add('_GoBack', $bookmarkRange)\n\n$sentinelRange = $d.Range($paragraphStart + $bodyLength, $paragraphStart + $bodyLength + $sentinel.Length)\n$sentinelRange.Delete()\n"}
